$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change style of A57 from "date only" to "date + time" style (matching previous rows)
$ws.Range("A57").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 58 with the new day's data, using the "date only" style previously on A57
$ws.Range("A58").Value = 45798
$ws.Range("A58").NumberFormat = "YYYY-MM-DD"
$ws.Range("B58").Value = 242
$ws.Range("C58").Value = 253
$ws.Range("D58").Value = 244
